$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 92.27273
$ws.Range("I2").Value = 90.14286
$ws.Range("J2").Value = 96
$ws.Range("K2").Value = 90.14286
$ws.Range("L2").Value = 96
$ws.Range("M2").Value = 22.85714
$ws.Range("N2").Value = -322
$ws.Range("H12").Value = 3667
$ws.Range("I12").Value = 5000
$ws.Range("K12").Value = 5000
$ws.Range("M12").Value = -4830
$ws.Range("H21").Value = 5017
$ws.Range("I21").Value = 5017
$ws.Range("K21").Value = 5017
$ws.Range("M21").Value = -4549
$ws.Range("H23").Value = 5017
$ws.Range("I23").Value = 5017
$ws.Range("K23").Value = 5017
$ws.Range("M23").Value = -4783
$ws.Range("H29").Value = 2000
$ws.Range("J29").Value = 3000
$ws.Range("L29").Value = 9000
$ws.Range("N29").Value = -9562
$ws.Range("H38").Value = 759.36365
$ws.Range("I38").Value = 510.6
$ws.Range("J38").Value = 966.6667
$ws.Range("K38").Value = 1531.8
$ws.Range("L38").Value = 2900.0001
$ws.Range("M38").Value = -1159.8
$ws.Range("N38").Value = -3644.0001
$ws.Range("H58").Value = 4216.25
$ws.Range("I58").Value = 243.33333
$ws.Range("J58").Value = 6600
$ws.Range("K58").Value = 729.99999
$ws.Range("L58").Value = 19800
$ws.Range("M58").Value = -579.99999
$ws.Range("N58").Value = -20100
$ws.Range("H87").Value = 43015.5
$ws.Range("J87").Value = 43015.5
$ws.Range("L87").Value = 43015.5
$ws.Range("N87").Value = -45511.5
$ws.Range("H90").Value = 43015.5
$ws.Range("J90").Value = 43015.5
$ws.Range("L90").Value = 129046.5
$ws.Range("N90").Value = -141526.5
$ws.Range("H125").Value = 2018
$ws.Range("I125").Value = 1000
$ws.Range("J125").Value = 3036
$ws.Range("K125").Value = 9000
$ws.Range("L125").Value = 27324
$ws.Range("M125").Value = -6540
$ws.Range("N125").Value = -32244
$ws.Range("H129").Value = 837.9722
$ws.Range("J129").Value = 837.9722
$ws.Range("L129").Value = 2513.9166
$ws.Range("N129").Value = -12513.9166
$ws.Range("H137").Value = 44828.914
$ws.Range("I137").Value = 1041.5
$ws.Range("K137").Value = 3124.5
$ws.Range("M137").Value = -574.5

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1511.92
$ws.Range("I122").Value = 1557.7273
$ws.Range("J122").Value = 1176
$ws.Range("K122").Value = 4673.1819
$ws.Range("L122").Value = 3528
$ws.Range("M122").Value = -2223.1819
$ws.Range("N122").Value = -8428

# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3967
$ws.Range("I20").Value = 5484.6665
$ws.Range("J20").Value = 2449.3333
$ws.Range("K20").Value = 5484.6665
$ws.Range("L20").Value = 2449.3333
$ws.Range("M20").Value = -5237.6665
$ws.Range("N20").Value = -2943.3333
$ws.Range("H81").Value = 22085.666
$ws.Range("J81").Value = 22085.666
$ws.Range("L81").Value = 22085.666
$ws.Range("N81").Value = -24207.666
$ws.Range("H84").Value = 22085.666
$ws.Range("J84").Value = 22085.666
$ws.Range("L84").Value = 66256.99800000001
$ws.Range("N84").Value = -76864.99800000001
$ws.Range("H107").Value = 797
$ws.Range("I107").Value = 695.5
$ws.Range("K107").Value = 695.5
$ws.Range("M107").Value = 1224.5

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1347.091
$ws.Range("I5").Value = 789.125
$ws.Range("J5").Value = 2835
$ws.Range("K5").Value = 2367.375
$ws.Range("L5").Value = 8505
$ws.Range("M5").Value = -2255.375
$ws.Range("N5").Value = -8729
$ws.Range("H14").Value = 296.18182
$ws.Range("I14").Value = 296.18182
$ws.Range("K14").Value = 888.54546
$ws.Range("M14").Value = -715.54546
$ws.Range("H23").Value = 572.2727
$ws.Range("J23").Value = 764.9286
$ws.Range("L23").Value = 2294.7858
$ws.Range("N23").Value = -2764.7858
$ws.Range("H26").Value = 290.54544
$ws.Range("I26").Value = 113.333336
$ws.Range("J26").Value = 503.2
$ws.Range("K26").Value = 340.000008
$ws.Range("L26").Value = 1509.6
$ws.Range("M26").Value = -52.00000799999998
$ws.Range("N26").Value = -2085.6
$ws.Range("H47").Value = 648.2
$ws.Range("I47").Value = 80.333336
$ws.Range("K47").Value = 241.000008
$ws.Range("M47").Value = 189.999992
$ws.Range("H87").Value = 15531.77
$ws.Range("I87").Value = 7676.625
$ws.Range("J87").Value = 28100
$ws.Range("K87").Value = 23029.875
$ws.Range("L87").Value = 84300
$ws.Range("M87").Value = -21781.875
$ws.Range("N87").Value = -86796
$ws.Range("H90").Value = 15531.77
$ws.Range("I90").Value = 7676.625
$ws.Range("J90").Value = 28100
$ws.Range("K90").Value = 69089.625
$ws.Range("L90").Value = 252900
$ws.Range("M90").Value = -62849.625
$ws.Range("N90").Value = -265380
$ws.Range("H114").Value = 522.4167
$ws.Range("I114").Value = 162.5
$ws.Range("J114").Value = 702.375
$ws.Range("K114").Value = 487.5
$ws.Range("L114").Value = 2107.125
$ws.Range("M114").Value = 2766.5
$ws.Range("N114").Value = -8615.125
$ws.Range("H117").Value = 2643
$ws.Range("I117").Value = 1464.5
$ws.Range("J117").Value = 5000
$ws.Range("K117").Value = 4393.5
$ws.Range("L117").Value = 15000
$ws.Range("M117").Value = -951.5
$ws.Range("N117").Value = -21884
$ws.Range("H118").Value = 125002280
$ws.Range("J118").Value = 4500
$ws.Range("L118").Value = 13500
$ws.Range("N118").Value = -15986
$ws.Range("H120").Value = 16671.666
$ws.Range("I120").Value = 10015
$ws.Range("K120").Value = 30045
$ws.Range("M120").Value = -25207
$ws.Range("H129").Value = 313840.94
$ws.Range("I129").Value = 800
$ws.Range("J129").Value = 358561.06
$ws.Range("K129").Value = 2400
$ws.Range("L129").Value = 1075683.18
$ws.Range("M129").Value = 2600
$ws.Range("N129").Value = -1085683.18

# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 2726226.5
$ws.Range("I70").Value = 12850.5
$ws.Range("K70").Value = 12850.5
$ws.Range("M70").Value = -12580.5
$ws.Range("H73").Value = 2726226.5
$ws.Range("I73").Value = 12850.5
$ws.Range("K73").Value = 12850.5
$ws.Range("M73").Value = -11914.5
$ws.Range("H131").Value = 682.77
$ws.Range("I131").Value = 586.75
$ws.Range("J131").Value = 686.7708
$ws.Range("K131").Value = 1760.25
$ws.Range("L131").Value = 2060.3124
$ws.Range("M131").Value = 3279.75
$ws.Range("N131").Value = -12140.3124
$ws.Range("H132").Value = 1106.2142
$ws.Range("I132").Value = 849.875
$ws.Range("J132").Value = 1448
$ws.Range("K132").Value = 7648.875
$ws.Range("L132").Value = 13032
$ws.Range("M132").Value = -5118.875
$ws.Range("N132").Value = -18092
$ws.Range("H135").Value = 1347.091
$ws.Range("I135").Value = 789.125
$ws.Range("J135").Value = 2835
$ws.Range("K135").Value = 7102.125
$ws.Range("L135").Value = 25515
$ws.Range("M135").Value = -4567.125
$ws.Range("N135").Value = -30585

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 351.7857
$ws.Range("I16").Value = 334
$ws.Range("J16").Value = 500
$ws.Range("K16").Value = 334
$ws.Range("L16").Value = 500
$ws.Range("M16").Value = -164
$ws.Range("N16").Value = -840
$ws.Range("H61").Value = 5870.6313
$ws.Range("I61").Value = 2282.5557
$ws.Range("K61").Value = 2282.5557
$ws.Range("M61").Value = -2080.5557
$ws.Range("H82").Value = 1705.6364
$ws.Range("I82").Value = 2019.6666
$ws.Range("J82").Value = 1488.2307
$ws.Range("K82").Value = 2019.6666
$ws.Range("L82").Value = 1488.2307
$ws.Range("M82").Value = -1658.6666
$ws.Range("N82").Value = -2210.2307
$ws.Range("H85").Value = 1705.6364
$ws.Range("I85").Value = 2019.6666
$ws.Range("J85").Value = 1488.2307
$ws.Range("K85").Value = 2019.6666
$ws.Range("L85").Value = 1488.2307
$ws.Range("M85").Value = -771.6666
$ws.Range("N85").Value = -3984.2307
$ws.Range("H113").Value = 5870.6313
$ws.Range("I113").Value = 2282.5557
$ws.Range("K113").Value = 2282.5557
$ws.Range("M113").Value = -112.5556999999999
$ws.Range("H136").Value = 84698.5
$ws.Range("I136").Value = 101398.2
$ws.Range("J136").Value = 1200
$ws.Range("K136").Value = 304194.6
$ws.Range("L136").Value = 3600
$ws.Range("M136").Value = -301644.6
$ws.Range("N136").Value = -8700

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()
